# Fruta / hortaliza, semanal
#
# The daily logic re-sorted the rows of this weekly subset: the values in
# columns D (Fecha), L (Calidad), M (Volumen), N (Precio mínimo),
# O (Precio máximo), P (Precio promedio ponderado), R (Origen) and
# S (Precio $/Kg) for rows 2..32 get reshuffled to a new row order
# (columns A, B, C, E, F, G, H, I, J, K, Q, T stay constant/unchanged).
#
# Snapshot the current ("before") values for the columns that move, then
# write them back out in the new ("after") order, so no data is lost
# even though several rows land on top of each other's old positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 32

# Columns that participate in the reshuffle.
$cols = @(4, 12, 13, 14, 15, 16, 18, 19)   # D, L, M, N, O, P, R, S

# destinationRow -> sourceRow (i.e. row $key ends up holding what row
# $value used to hold, for the columns listed above)
$rowMap = @{
    2 = 24; 3 = 25; 4 = 19; 5 = 20; 6 = 23; 7 = 17; 8 = 18; 9 = 22; 10 = 26;
    11 = 27; 12 = 28; 13 = 32; 14 = 13; 15 = 14; 16 = 15; 17 = 12; 18 = 16;
    19 = 8; 20 = 9; 21 = 7; 22 = 29; 23 = 5; 24 = 6; 25 = 2; 26 = 3; 27 = 21;
    28 = 4; 29 = 10; 30 = 11; 31 = 30; 32 = 31
}

# 1) Snapshot every source cell's current value before anything is written.
#    (Use Value2 - plain .Value has been observed to mis-resolve on read
#    in this engine, returning the property descriptor instead of the
#    underlying data; Value2 reads/writes both numbers and strings fine.)
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# 2) Write the snapshotted values back into their new row positions.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $rowMap[$r]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $snapshot["$src,$c"]
    }
}
